$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "model_8_6_24"
$ws.Range("B2").Value = -0.09273789730344206
$ws.Range("C2").Value = -1.012154730746021
$ws.Range("D2").Value = -4.152979013636195
$ws.Range("E2").Value = -2.10887844408628
$ws.Range("F2").Value = 1.209338903427124
$ws.Range("G2").Value = 2.934303045272827
$ws.Range("H2").Value = 5.164154052734375
$ws.Range("I2").Value = 3.983642101287842
$ws.Range("A3").Value = "model_8_6_23"
$ws.Range("B3").Value = -0.09030335439939163
$ws.Range("C3").Value = -1.015784341513713
$ws.Range("D3").Value = -4.128178747260187
$ws.Range("E3").Value = -2.101941350584962
$ws.Range("F3").Value = 1.206644654273987
$ws.Range("G3").Value = 2.939595937728882
$ws.Range("H3").Value = 5.139300346374512
$ws.Range("I3").Value = 3.974752426147461
$ws.Range("A4").Value = "model_8_6_22"
$ws.Range("B4").Value = -0.08470543407412667
$ws.Range("C4").Value = -1.018818555343681
$ws.Range("D4").Value = -4.079915590246125
$ws.Range("E4").Value = -2.08600168292041
$ws.Range("F4").Value = 1.200449228286743
$ws.Range("G4").Value = 2.944020509719849
$ws.Range("H4").Value = 5.090932369232178
$ws.Range("I4").Value = 3.954327821731567
$ws.Range("A5").Value = "model_8_6_21"
$ws.Range("B5").Value = -0.07769226302917875
$ws.Range("C5").Value = -1.020563179263584
$ws.Range("D5").Value = -4.022793170544749
$ws.Range("E5").Value = -2.066032138658124
$ws.Range("F5").Value = 1.192687749862671
$ws.Range("G5").Value = 2.946564912796021
$ws.Range("H5").Value = 5.03368616104126
$ws.Range("I5").Value = 3.928739309310913
$ws.Range("A6").Value = "model_8_6_20"
$ws.Range("B6").Value = -0.06387007243716614
$ws.Range("C6").Value = -1.015575081433473
$ws.Range("D6").Value = -3.924073054507711
$ws.Range("E6").Value = -2.026690612171861
$ws.Range("F6").Value = 1.177390694618225
$ws.Range("G6").Value = 2.939290761947632
$ws.Range("H6").Value = 4.934751510620117
$ws.Range("I6").Value = 3.878328561782837
$ws.Range("A7").Value = "model_8_6_19"
$ws.Range("B7").Value = -0.05315664012430332
$ws.Range("C7").Value = -1.009541546640976
$ws.Range("D7").Value = -3.85106836344803
$ws.Range("E7").Value = -1.996187083326114
$ws.Range("F7").Value = 1.165534138679504
$ws.Range("G7").Value = 2.930492162704468
$ws.Range("H7").Value = 4.861588478088379
$ws.Range("I7").Value = 3.839241981506348
$ws.Range("A8").Value = "model_8_6_18"
$ws.Range("B8").Value = -0.04435932903006545
$ws.Range("C8").Value = -1.014593598327592
$ws.Range("D8").Value = -3.774704730656787
$ws.Range("E8").Value = -1.971129954731939
$ws.Range("F8").Value = 1.155798077583313
$ws.Range("G8").Value = 2.937859535217285
$ws.Range("H8").Value = 4.785059452056885
$ws.Range("I8").Value = 3.807134389877319
$ws.Range("A9").Value = "model_8_6_17"
$ws.Range("B9").Value = -0.0358482577767425
$ws.Range("C9").Value = -1.011144779179304
$ws.Range("D9").Value = -3.714480166878709
$ws.Range("E9").Value = -1.94688395044356
$ws.Range("F9").Value = 1.146378755569458
$ws.Range("G9").Value = 2.932830095291138
$ws.Range("H9").Value = 4.724704265594482
$ws.Range("I9").Value = 3.776065587997437
$ws.Range("A10").Value = "model_8_6_16"
$ws.Range("B10").Value = -0.01090532409447187
$ws.Range("C10").Value = -0.9942548280355656
$ws.Range("D10").Value = -3.54922770778419
$ws.Range("E10").Value = -1.875887608329849
$ws.Range("F10").Value = 1.11877429485321
$ws.Range("G10").Value = 2.908199787139893
$ws.Range("H10").Value = 4.559093475341797
$ws.Range("I10").Value = 3.685092926025391
$ws.Range("A11").Value = "model_8_6_15"
$ws.Range("B11").Value = 0.000145806180597341
$ws.Range("C11").Value = -0.9972879768962353
$ws.Range("D11").Value = -3.458705966164477
$ws.Range("E11").Value = -1.844398629655291
$ws.Range("F11").Value = 1.106543898582458
$ws.Range("G11").Value = 2.912622928619385
$ws.Range("H11").Value = 4.468375205993652
$ws.Range("I11").Value = 3.6447434425354
$ws.Range("A12").Value = "model_8_6_1"
$ws.Range("B12").Value = 0.009097935982429139
$ws.Range("C12").Value = -2.321536677893954
$ws.Range("D12").Value = -1.159216566688853
$ws.Range("E12").Value = -1.795939433482639
$ws.Range("F12").Value = 1.096636652946472
$ws.Range("G12").Value = 4.843760013580322
$ws.Range("H12").Value = 2.163899183273315
$ws.Range("I12").Value = 3.582648992538452
$ws.Range("A13").Value = "model_8_6_14"
$ws.Range("B13").Value = 0.01528585696714746
$ws.Range("C13").Value = -1.000824770065175
$ws.Range("D13").Value = -3.33572192794866
$ws.Range("E13").Value = -1.801266797553361
$ws.Range("F13").Value = 1.089788317680359
$ws.Range("G13").Value = 2.917780637741089
$ws.Range("H13").Value = 4.3451247215271
$ws.Range("I13").Value = 3.589475631713867
$ws.Range("A14").Value = "model_8_6_13"
$ws.Range("B14").Value = 0.03246333944382696
$ws.Range("C14").Value = -1.003513234407131
$ws.Range("D14").Value = -3.198345661526156
$ws.Range("E14").Value = -1.752322869228246
$ws.Range("F14").Value = 1.070778012275696
$ws.Range("G14").Value = 2.921701192855835
$ws.Range("H14").Value = 4.207450866699219
$ws.Range("I14").Value = 3.52675986289978
$ws.Range("A15").Value = "model_8_6_12"
$ws.Range("B15").Value = 0.0811863129757705
$ws.Range("C15").Value = -0.9937622171911809
$ws.Range("D15").Value = -2.837418631150658
$ws.Range("E15").Value = -1.613611781538211
$ws.Range("F15").Value = 1.016856074333191
$ws.Range("G15").Value = 2.90748119354248
$ws.Range("H15").Value = 3.845740795135498
$ws.Range("I15").Value = 3.349018812179565
$ws.Range("A16").Value = "model_8_6_11"
$ws.Range("B16").Value = 0.1046347754240831
$ws.Range("C16").Value = -0.9948778727992018
$ws.Range("D16").Value = -2.654031571793501
$ws.Range("E16").Value = -1.54678700014497
$ws.Range("F16").Value = 0.9909054636955261
$ws.Range("G16").Value = 2.909108400344849
$ws.Range("H16").Value = 3.661956071853638
$ws.Range("I16").Value = 3.263391494750977
$ws.Range("A17").Value = "model_8_6_10"
$ws.Range("B17").Value = 0.124155038623934
$ws.Range("C17").Value = -0.9948211829593883
$ws.Range("D17").Value = -2.50285498282916
$ws.Range("E17").Value = -1.49111220398112
$ws.Range("F17").Value = 0.969302237033844
$ws.Range("G17").Value = 2.9090256690979
$ws.Range("H17").Value = 3.510451555252075
$ws.Range("I17").Value = 3.19205117225647
$ws.Range("A18").Value = "model_8_6_9"
$ws.Range("B18").Value = 0.1438363840560252
$ws.Range("C18").Value = -0.9949434837435545
$ws.Range("D18").Value = -2.350051151162876
$ws.Range("E18").Value = -1.434948407471282
$ws.Range("F18").Value = 0.9475207924842834
$ws.Range("G18").Value = 2.909203767776489
$ws.Range("H18").Value = 3.357316255569458
$ws.Range("I18").Value = 3.120084285736084
$ws.Range("A19").Value = "model_8_6_8"
$ws.Range("B19").Value = 0.1665770890757999
$ws.Range("C19").Value = -0.9923662916159066
$ws.Range("D19").Value = -2.177904756252021
$ws.Range("E19").Value = -1.370034590913253
$ws.Range("F19").Value = 0.9223536252975464
$ws.Range("G19").Value = 2.90544581413269
$ws.Range("H19").Value = 3.184796810150146
$ws.Range("I19").Value = 3.036905288696289
$ws.Range("A20").Value = "model_8_6_0"
$ws.Range("B20").Value = 0.1853199283319504
$ws.Range("C20").Value = -2.015074083680278
$ws.Range("D20").Value = -0.1478906653596841
$ws.Range("E20").Value = -1.239075692634021
$ws.Range("F20").Value = 0.901610791683197
$ws.Range("G20").Value = 4.396848678588867
$ws.Range("H20").Value = 1.15038013458252
$ws.Range("I20").Value = 2.869097232818604
$ws.Range("A21").Value = "model_8_6_7"
$ws.Range("B21").Value = 0.1887480847695719
$ws.Range("C21").Value = -0.9842457378743343
$ws.Range("D21").Value = -2.019082006905151
$ws.Range("E21").Value = -1.306686484853
$ws.Range("F21").Value = 0.8978167772293091
$ws.Range("G21").Value = 2.893603324890137
$ws.Range("H21").Value = 3.02562952041626
$ws.Range("I21").Value = 2.955732345581055
$ws.Range("A22").Value = "model_8_6_6"
$ws.Range("B22").Value = 0.2077653619928833
$ws.Range("C22").Value = -0.9696787661746324
$ws.Range("D22").Value = -1.895094509875312
$ws.Range("E22").Value = -1.252275903321601
$ws.Range("F22").Value = 0.8767703175544739
$ws.Range("G22").Value = 2.872360706329346
$ws.Range("H22").Value = 2.901372909545898
$ws.Range("I22").Value = 2.886011838912964
$ws.Range("A23").Value = "model_8_6_5"
$ws.Range("B23").Value = 0.2214652175370833
$ws.Range("C23").Value = -0.9649611959164845
$ws.Range("D23").Value = -1.795802959574566
$ws.Range("E23").Value = -1.212891035190357
$ws.Range("F23").Value = 0.861608624458313
$ws.Range("G23").Value = 2.86548113822937
$ws.Range("H23").Value = 2.801866292953491
$ws.Range("I23").Value = 2.835545063018799
$ws.Range("A24").Value = "model_8_6_4"
$ws.Range("B24").Value = 0.2396098266747348
$ws.Range("C24").Value = -0.9524351398878377
$ws.Range("D24").Value = -1.674569096874512
$ws.Range("E24").Value = -1.160723943168089
$ws.Range("F24").Value = 0.8415278792381287
$ws.Range("G24").Value = 2.847214698791504
$ws.Range("H24").Value = 2.68036937713623
$ws.Range("I24").Value = 2.768699407577515
$ws.Range("A25").Value = "model_8_6_3"
$ws.Range("B25").Value = 0.2557104871934369
$ws.Range("C25").Value = -0.9380508757161479
$ws.Range("D25").Value = -1.571980161711036
$ws.Range("E25").Value = -1.11430152862168
$ws.Range("F25").Value = 0.8237091898918152
$ws.Range("G25").Value = 2.826237916946411
$ws.Range("H25").Value = 2.577558040618896
$ws.Range("I25").Value = 2.709214925765991
$ws.Range("A26").Value = "model_8_6_2"
$ws.Range("B26").Value = 0.3739783326781468
$ws.Range("C26").Value = -0.7190988171251351
$ws.Range("D26").Value = -1.010588277514764
$ws.Range("E26").Value = -0.7757617423268606
$ws.Range("F26").Value = 0.6928215026855469
$ws.Range("G26").Value = 2.506942749023438
$ws.Range("H26").Value = 2.014948844909668
$ws.Range("I26").Value = 2.275418043136597
